$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4360
$ws.Range("J43").Value = 5280
$ws.Range("L43").Value = 5280
$ws.Range("N43").Value = -5418
$ws.Range("H86").Value = 57978850
$ws.Range("I86").Value = 20841250
$ws.Range("J86").Value = 142864780
$ws.Range("K86").Value = 20841250
$ws.Range("L86").Value = 142864780
$ws.Range("M86").Value = -20840127
$ws.Range("N86").Value = -142867026
$ws.Range("H89").Value = 57978850
$ws.Range("I89").Value = 20841250
$ws.Range("J89").Value = 142864780
$ws.Range("K89").Value = 104206250
$ws.Range("L89").Value = 714323900
$ws.Range("M89").Value = -104200634
$ws.Range("N89").Value = -714335132
$ws.Range("H106").Value = 373026.28
$ws.Range("I106").Value = 506617.9
$ws.Range("J106").Value = 5649.25
$ws.Range("K106").Value = 506617.9
$ws.Range("L106").Value = 5649.25
$ws.Range("M106").Value = -505986.9
$ws.Range("N106").Value = -6911.25
$ws.Range("H118").Value = 733.1667
$ws.Range("I118").Value = 742.4
$ws.Range("K118").Value = 2227.2
$ws.Range("M118").Value = -570.1999999999998
$ws.Range("H132").Value = 2005.1852
$ws.Range("J132").Value = 4999.25
$ws.Range("L132").Value = 14997.75
$ws.Range("N132").Value = -20057.75
$ws.Range("H136").Value = 96491.664
$ws.Range("J136").Value = 96491.664
$ws.Range("L136").Value = 96491.664
$ws.Range("N136").Value = -106691.664
$ws.Range("H138").Value = 2228.889
$ws.Range("J138").Value = 2993.353
$ws.Range("L138").Value = 8980.059000000001
$ws.Range("N138").Value = -19260.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12199.478
$ws.Range("I32").Value = 7017.8237
$ws.Range("J32").Value = 29817.1
$ws.Range("K32").Value = 7017.8237
$ws.Range("L32").Value = 29817.1
$ws.Range("M32").Value = -6730.8237
$ws.Range("N32").Value = -30391.1
$ws.Range("H132").Value = 1595.5358
$ws.Range("I132").Value = 1434.48
$ws.Range("K132").Value = 4303.440000000001
$ws.Range("M132").Value = -1773.440000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 557.5
$ws.Range("I94").Value = 550.9167
$ws.Range("K94").Value = 550.9167
$ws.Range("M94").Value = -99.91669999999999
$ws.Range("H140").Value = 99847
$ws.Range("J140").Value = 99847
$ws.Range("L140").Value = 99847
$ws.Range("N140").Value = -110207

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12121.588
$ws.Range("I31").Value = 1912.88
$ws.Range("J31").Value = 40479.11
$ws.Range("K31").Value = 1912.88
$ws.Range("L31").Value = 40479.11
$ws.Range("M31").Value = -1617.88
$ws.Range("N31").Value = -41069.11
$ws.Range("H34").Value = 12121.588
$ws.Range("I34").Value = 1912.88
$ws.Range("J34").Value = 40479.11
$ws.Range("K34").Value = 1912.88
$ws.Range("L34").Value = 40479.11
$ws.Range("M34").Value = -1710.88
$ws.Range("N34").Value = -40883.11
$ws.Range("H62").Value = 4449.5
$ws.Range("I62").Value = 3899
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3899
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3275
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4449.5
$ws.Range("I65").Value = 3899
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 19495
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16375
$ws.Range("N65").Value = -31240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4506.8
$ws.Range("I34").Value = 169
$ws.Range("J34").Value = 7398.6665
$ws.Range("K34").Value = 507
$ws.Range("L34").Value = 22195.9995
$ws.Range("M34").Value = -423
$ws.Range("N34").Value = -22363.9995
$ws.Range("H38").Value = 27.857143
$ws.Range("I38").Value = 16.666666
$ws.Range("J38").Value = 42.77778
$ws.Range("K38").Value = 49.999998
$ws.Range("L38").Value = 128.33334
$ws.Range("M38").Value = 297.000002
$ws.Range("N38").Value = -822.33334
$ws.Range("H131").Value = 33929.645
$ws.Range("I131").Value = 500499.5
$ws.Range("K131").Value = 1501498.5
$ws.Range("M131").Value = -1496458.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 20011000
$ws.Range("I20").Value = 999.5
$ws.Range("J20").Value = 40021000
$ws.Range("K20").Value = 999.5
$ws.Range("L20").Value = 40021000
$ws.Range("M20").Value = -754.5
$ws.Range("N20").Value = -40021490
$ws.Range("H24").Value = 26000
$ws.Range("I24").Value = 26000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 26000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -25827
$ws.Range("N24").ClearContents()
$ws.Range("H97").Value = 527291.5600000001
$ws.Range("I97").Value = 769536
$ws.Range("J97").Value = 2428.6667
$ws.Range("K97").Value = 769536
$ws.Range("L97").Value = 2428.6667
$ws.Range("M97").Value = -769040
$ws.Range("N97").Value = -3420.6667
$ws.Range("H113").Value = 4764106.5
$ws.Range("I113").Value = 1199.5
$ws.Range("J113").Value = 6669269
$ws.Range("K113").Value = 1199.5
$ws.Range("L113").Value = 6669269
$ws.Range("M113").Value = 970.5
$ws.Range("N113").Value = -6673609

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3327.7
$ws.Range("I16").Value = 3141.889
$ws.Range("K16").Value = 3141.889
$ws.Range("M16").Value = -2971.889
$ws.Range("H22").Value = 987.93335
$ws.Range("I22").Value = 849.93335
$ws.Range("K22").Value = 849.93335
$ws.Range("M22").Value = -554.93335
$ws.Range("H27").Value = 987.93335
$ws.Range("I27").Value = 849.93335
$ws.Range("K27").Value = 849.93335
$ws.Range("M27").Value = -742.93335
$ws.Range("H46").Value = 2499.75
$ws.Range("I46").Value = 1499.5
$ws.Range("J46").Value = 3500
$ws.Range("K46").Value = 1499.5
$ws.Range("L46").Value = 3500
$ws.Range("M46").Value = -1311.5
$ws.Range("N46").Value = -3876
$ws.Range("H55").Value = 1730.1428
$ws.Range("I55").Value = 929.5172
$ws.Range("K55").Value = 929.5172
$ws.Range("M55").Value = -756.5172
$ws.Range("H100").Value = 75999.8
$ws.Range("I100").Value = 93749.75
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 93749.75
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -93208.75
$ws.Range("N100").Value = -6082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
